$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 44581.62222222222
$ws.Range("N2").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"

$ws.Range("N3").Value = 44581
$ws.Range("N3").NumberFormat = "mm-dd-yy"

$ws.Columns.Item(14).ColumnWidth = 15

$ws.Range("N3").Select()
$ws.Application.ActiveWindow.ScrollColumn = 7
